# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: the data table currently spans rows 16-25 (10 rows,
#    last one carrying the special "closing border" style). The new data
#    has 14 rows, so insert 4 blank rows right before the closing row
#    (row 25) -- this pushes row 25 (and everything below it, incl. the
#    merged footer block) down to row 29, exactly like Excel would when a
#    user inserts rows above the last line of a bordered table.
$ws.Rows("25:28").Insert()

# Copy the normal-row formatting (from row 24, a "middle" row of the old
# table) onto the freshly inserted rows 25-28 so every column keeps the
# right borders/number formats.
$ws.Range("B24:J24").Copy()
$ws.Range("B25:J28").PasteSpecial(-4122)

# --- 2. Write the new data set (rows 16-29) ---
$data = @(
    @("CC","45781333","JACKELIN BARRIOS YEPEZ","2104",15748,1160000),
    @("CC","45781333","JACKELIN BARRIOS YEPEZ","2103",36341,1160000),
    @("CC","33273074","KATIA MERCEDES CONTRERAS ARDILA","2011",35112,877803),
    @("CC","33337424","MARBEL LUZ BALLESTAS BUELVAS","2011",35112,908526),
    @("CC","33337424","MARBEL LUZ BALLESTAS BUELVAS","2010",35112,908526),
    @("CC","33341705","DIANORA ELENA MARTINEZ MEZA","2104",15748,908526),
    @("CC","33341705","DIANORA ELENA MARTINEZ MEZA","2103",36341,908526),
    @("CC","1049942967","ANAYIBIS PEREZ HERNANDEZ","2011",10534,877803),
    @("CC","1049942967","ANAYIBIS PEREZ HERNANDEZ","2010",35112,877803),
    @("CC","1049939325","MARIA VANESSA BENAVIDES MARIMON","2011",10534,877803),
    @("CC","1049939325","MARIA VANESSA BENAVIDES MARIMON","2010",35112,877803),
    @("CC","41371858","MARIA TERESA PADILLA CAÃ?ATE","2102",1211,908526),
    @("CC","41371858","MARIA TERESA PADILLA CAÃ?ATE","2011",35112,908526),
    @("CC","41371858","MARIA TERESA PADILLA CAÃ?ATE","2010",35112,908526)
)

$r = 16
foreach ($row in $data) {
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $r = $r + 1
}

# --- 3. Update the summary cells above the table ---
# Valor Mora total (E11) = sum of Valor Mora column
$ws.Range("E11").Value = 372241
# Cant. Trabajadores (C13) / Cant. Periodos (F13)
$ws.Range("C13").Value = 7
$ws.Range("F13").Value = 5

Write-Output "edit applied"
